# Refactorizacion de estructura de scrapping
# Update job-listing data (rows 2-5) with a refreshed scrape, and adjust
# header/row heights to fit the new (sometimes multi-line) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Practicante desarrollador de software web (ROWAN NETWORKS)
$ws.Range("A2").Value = "Practicante desarrollador de software web - job post"
$ws.Range("B2").Value = "ROWAN NETWORKS"
$ws.Range("C2").Value = "`$8000 por mes - TiempoCompleto"

# Row 3 - Desarrollador/a (FINANCIERA CRECE CON VALE)
$ws.Range("A3").Value = "Desarrollador/a"
$ws.Range("B3").Value = "FINANCIERA CRECE CON VALE"
$ws.Range("C3").Value = "Tiempo completo"

# Row 4 - Desarrollador de Software - job post (MOVERET INCORPORATED)
$ws.Range("A4").Value = "Desarrollador de Software`n- job post"
$ws.Range("B4").Value = "MOVERET INCORPORATED"
$ws.Range("C4").Value = "Desde `$700,000 por año - Tiempo completo"

# Row 5 - Desarrollador web - job post (SISTEMAS DE SEGURIDAD PRIVADA ARGOS S.A. DE C.V.)
$ws.Range("A5").Value = "Desarrollador web`n- job post"
$ws.Range("B5").Value = "SISTEMAS DE SEGURIDAD PRIVADA ARGOS S.A. DE C.V."
$ws.Range("C5").Value = "`$13,000 a `$15,000 por mes - Tiempo completo"

# Row heights adjustment
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 31.5
$ws.Rows.Item(3).RowHeight = 31.5
$ws.Rows.Item(4).RowHeight = 31.5
$ws.Rows.Item(5).RowHeight = 31.5
